$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1785
$ws.Range("F3").Value = 10438
$ws.Range("F5").Value = 18
$ws.Range("F6").Value = 617
$ws.Range("F8").Value = 1737
$ws.Range("F9").Value = 441
$ws.Range("F10").Value = 13
$ws.Range("F11").Value = 253
$ws.Range("F12").Value = 545
$ws.Range("F13").Value = 1183
$ws.Range("F14").Value = 150
$ws.Range("F15").Value = 37
$ws.Range("F16").Value = 1031
$ws.Range("F17").Value = 36
$ws.Range("F18").Value = 123
$ws.Range("F19").Value = 428
$ws.Range("F20").Value = 428
$ws.Range("F22").Value = 362
$ws.Range("F23").Value = 62
$ws.Range("F24").Value = 1099
$ws.Range("F25").Value = 1125
$ws.Range("F26").Value = 1230
$ws.Range("F27").Value = 219
$ws.Range("F28").Value = 1426
$ws.Range("F29").Value = 725
$ws.Range("F30").Value = 263
$ws.Range("F31").Value = 33
$ws.Range("F33").Value = 701
$ws.Range("F34").Value = 263
$ws.Range("F35").Value = 745
$ws.Range("F37").Value = 798
$ws.Range("F38").Value = 135595
$ws.Range("F39").Value = 814
$ws.Range("F40").Value = 534
$ws.Range("F44").Value = 1392

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 96
$ws.Range("F11").Value = 75
$ws.Range("F14").Value = 119
$ws.Range("F19").Value = 1157
$ws.Range("F21").Value = 2293
$ws.Range("F23").Value = 354
$ws.Range("F24").Value = 697
$ws.Range("F30").Value = 385
$ws.Range("F33").Value = 229
$ws.Range("F43").Value = 15
$ws.Range("F46").Value = 91

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 833
$ws.Range("F5").Value = 215
$ws.Range("F6").Value = 2582
$ws.Range("F7").Value = 4270
$ws.Range("F8").Value = 81
$ws.Range("F10").Value = 421
$ws.Range("F11").Value = 387
$ws.Range("F12").Value = 278
$ws.Range("F13").Value = 192
$ws.Range("F14").Value = 90

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1785
$ws.Range("F3").Value = 833
$ws.Range("F4").Value = 10438
$ws.Range("F5").Value = 215
$ws.Range("F6").Value = 4270
$ws.Range("F7").Value = 617
$ws.Range("F8").Value = 387
$ws.Range("F9").Value = 1737
$ws.Range("F10").Value = 253
$ws.Range("F12").Value = 75
$ws.Range("F13").Value = 192
$ws.Range("F14").Value = 90
$ws.Range("F15").Value = 119
$ws.Range("F16").Value = 1031
$ws.Range("F17").Value = 36
$ws.Range("F18").Value = 123
$ws.Range("F19").Value = 428
$ws.Range("F20").Value = 428
$ws.Range("F21").Value = 362
$ws.Range("F22").Value = 2293
$ws.Range("F23").Value = 2293
$ws.Range("F25").Value = 1099
$ws.Range("F26").Value = 1125
$ws.Range("F27").Value = 1230
$ws.Range("F30").Value = 1426
$ws.Range("F31").Value = 725
$ws.Range("F32").Value = 385
$ws.Range("F33").Value = 701
$ws.Range("F35").Value = 745
$ws.Range("F37").Value = 798
$ws.Range("F38").Value = 229
$ws.Range("F39").Value = 814
$ws.Range("F40").Value = 534
$ws.Range("F44").Value = 1393
$ws.Range("F52").Value = 91
